$d = $word.ActiveDocument

# Paragraph 1: Title - "Questions: Rationalizing the denominator"
$r1 = $d.Paragraphs.Item(1).Range
$r1.Find.Execute("Questions: Rationalizing the denominator", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: Rationalizing the denominator", 2)

# Paragraph 2: Author - "Maximilian Volmar"
$r2 = $d.Paragraphs.Item(2).Range
$r2.Find.Execute("Maximilian Volmar", $false, $false, $false, $false, $false, $true, 1, $false, "Maximilian Volmar", 2)

# Paragraph 4: Abstract - "A selection of questions for the study guide on rationalizing the denominator."
$r4 = $d.Paragraphs.Item(4).Range
$r4.Find.Execute("A selection of questions for the study guide on rationalizing the denominator.", $false, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on rationalizing the denominator.", 2)

Write-Host "done"
